$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a 4th kernel ("Sigmoid") column of SVM results for the "c = 0.1" block
# (rows 27-50 of the worksheet), mirroring the Linear/Polynomial/RBF columns
# already present in columns B/C/D.
# ---------------------------------------------------------------------------

# 1vs2
$ws.Range("E28").Value = "SVMAccuracy : 60.6061"
$ws.Range("E29").Value = "SVMConfusionMatrix"
$ws.Range("E30").Value = "73  27"
$ws.Range("E31").Value = "52  48"

# 1vs3
$ws.Range("E34").Value = "SVMAccuracy : 53.0303"
$ws.Range("E35").Value = "SVMConfusionMatrix"
$ws.Range("E36").Value = "67  33"
$ws.Range("E37").Value = "61  39"

# 2vs3
$ws.Range("E40").Value = "SVMAccuracy : 51.5152"
$ws.Range("E41").Value = "SVMConfusionMatrix"
$ws.Range("E42").Value = "52  48"
$ws.Range("E43").Value = "48  52"

# 1vs2vs3
$ws.Range("E46").Value = "SVMAccuracy : 34.3434"
$ws.Range("E47").Value = "SVMConfusionMatrix"
$ws.Range("E48").Value = "21  48  30"
$ws.Range("E49").Value = "42  42  15"
$ws.Range("E50").Value = "24  36  39"

# New column E width to roughly match the other data columns
$ws.Columns.Item(5).ColumnWidth = 19.2

# View state: scroll the window up and move the active selection
$ws.Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("F39").Select()
